# Add new character styles
$d = $word.ActiveDocument

$styleGaN = $d.Styles.Add("GaNStyle", 2)
$styleGaN.Font.Name = "Calibri"
$styleGaN.Font.Size = 14

$styleParagraph = $d.Styles.Add("GaNParagraph", 2)
$styleParagraph.Font.Name = "Calibri"
$styleParagraph.Font.Size = 10

$styleLinks = $d.Styles.Add("GaNLinks", 2)
$styleLinks.Font.Name = "Calibri"
$styleLinks.Font.Bold = $true
$styleLinks.Font.Color = 8388608
$styleLinks.Font.Size = 9.5
$styleLinks.Font.Underline = 1

# Apply GaNStyle to all occurrences of the "V roku 2022..." paragraph text
$searchGaN = "V roku 2022 môžete pozorovať Súhvezdie Perzeus: 16. – 25. januára, 7. – 16. novembra, 6. – 15. decembra"
$range = $d.Content
$found = $range.Find.Execute($searchGaN, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $found = $range.Find.Execute($searchGaN, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# Apply GaNParagraph to the "Stávate sa súčasťou..." paragraph text
$searchParagraph = "Stávate sa súčasťou celosvetovej kampane Globe at Night, ktorej cieľom je meranie svetelného znečistenia. Pozorovaním  Súhvezdie Perzeus na nočnej oblohe a porovnávaním skutočnej situácie s našimi mapkami sa nielenže dozviete, ako osvetlenie vo Vašom okolí prispieva k svetelnému znečisteniu, ale budete môcť porovnať úroveň svetelného znečistenia aj s inými lokalitami z celého sveta. Vaše pozorovanie tiež rozšíri online databázu dokumentujúcu viditeľnosť nočnej oblohy na našej planéte"
$range2 = $d.Content
$found2 = $range2.Find.Execute($searchParagraph, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found2) {
    $range2.Style = "GaNParagraph"
    $range2.Collapse(0)
    $found2 = $range2.Find.Execute($searchParagraph, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# Apply GaNLinks to the "Mapky v tomto dokumente..." paragraph text
$searchLinks = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$found3 = $range3.Find.Execute($searchLinks, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found3) {
    $range3.Style = "GaNLinks"
    $range3.Collapse(0)
    $found3 = $range3.Find.Execute($searchLinks, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

Write-Host "Done applying GaN styles"
